# Apply "repull data, push all data, mean calculation" update.
# This re-pulls the dSF (column F) values for a subset of rows to reflect
# newly recomputed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new dSF (column F) value
$updates = @{
    2  = -1
    4  = -6
    5  = -1
    6  = -3
    8  = 5
    9  = -4
    12 = -2
    14 = -2
    23 = 4
    25 = -3
    27 = 3
    28 = 0
    32 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
